# Edit script: "Edit Linear Regression data for tech Industries and Add Logistic Regression in R file"
#
# Changes:
#  1. RawData!B2:B206 and RawData!R2:R206 get replaced with new raw data values
#     (R2 changes from a shared-string "NA" to a numeric value).
#  2. Selection/active-cell + active-sheet bookkeeping:
#     - Data sheet selection -> E10
#     - RawData sheet becomes the active/selected tab, selection -> D5
#     - HOWTO sheet loses tabSelected
#  3. Unused "NA" shared string gets pruned automatically once R2 no longer
#     references it.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsRaw  = $wb.Worksheets.Item("RawData")
$wsHow  = $wb.Worksheets.Item("HOWTO")

# ---------------------------------------------------------------------------
# 1. Bulk-replace the RawData B (tech-sector linear regression X) and
#    R (S&P500 logistic regression Y) columns with the new raw values.
# ---------------------------------------------------------------------------
$bCsv = "120.2813,136.8125,116.0625,106.6875,81.6875,74.5625,87.4,62,53.3,63.7,60.4,61.5,57,49.5,38.65,47,54,53.2,53.7,45.9,49.55,43.45,41.76,36.049999999999997,32.700000000000003,32.53,26.6,32.25,38.08,32.5,32.14,32.56,32.49,35.26,38.89,38.61,40.85,43.49,43.13,47.3,48.08,48.51,50.5,48.94,47.45,44.46,47.13,48.2,44.02,41.71,42.97,45.44,47.95,48.46,45.61,45.86,44.8,42.9,46.81,45.77,48.67,48.33,48.64,47.46,51,49.72,51.79,51.14,52.260100000000001,51.94,47.99,47.46,45.74,49.39,51.35,53.17,55.22,54.45,55.33,53.82,53.95,56.62,59.04,59.52,58.89,60.57,62.8,67.09,61.75,62.46,54.38,52.4,52.5,56.43,59.96,54.19,53.49,54.8,47.43,39.49,34.549999999999997,35.33,34.35,33.15,37.03,41.61,42.69,44.49,48.51,49.62,51.95,51.47,54.12,57.54,52.49,55.07,58.39,59.56,55.13,51.6,55.44,51.54,57.85,61.67,60.94,64.38,66.930000000000007,68.31,66.540000000000006,68.400000000000006,67.180000000000007,65.02,65.73,61.36,59.14,66.14,64.78,63.9,69.19,74.17,77.81,75.709999999999994,69.650000000000006,71.23,71.52,75.39,75.87,70.290000000000006,70.86,70.72,71.73,71.849999999999994,73.31,73.63,76.77,73.63,77.25,76.77,78.8,82.281999999999996,85.2,88.44,86.93,90.974000000000004,91.06,90.88,94.24,96.88,97.85,102,100.76,101.63,106.83,104.4,100.36,108.77,104.96,107.27,109.78,104.67,106.91,100.64,98.9,109.24,110.3,107.03,101.01,99.74,108.5,102.28,108.09,105.28,114.13,116.57,119,118.59,118.94,120.25,125.92,132.15,135.16999999999999,138.21,144.19999999999999,139.72,144.81"
$rCsv = "1430.83,1517.68,1436.51,1429.4,1314.95,1320.28,1366.01,1239.94,1160.33,1249.46,1255.82,1224.42,1211.23,1133.58,1040.94,1059.78,1139.45,1148.08,1130.2,1106.73,1147.3900000000001,1076.92,1067.1400000000001,989.81,911.62,916.07,815.28,885.76,936.31,879.82,855.7,841.15,848.18,916.92,963.59,974.5,990.31,1008.01,995.97,1050.71,1058.2,1111.92,1131.1300000000001,1144.94,1126.21,1107.3,1120.68,1140.8399999999999,1101.72,1104.24,1114.58,1130.2,1173.82,1211.92,1181.27,1203.5999999999999,1180.5899999999999,1156.8499999999999,1191.5,1191.33,1234.18,1220.33,1228.81,1207.01,1249.48,1248.29,1280.08,1280.6600000000001,1294.83,1310.6099999999999,1270.0899999999999,1270.2,1276.6600000000001,1303.82,1335.85,1377.94,1400.63,1418.3,1438.24,1406.82,1420.86,1482.37,1530.62,1503.35,1455.27,1473.99,1526.75,1549.38,1481.14,1468.36,1378.55,1330.63,1322.7,1385.59,1400.38,1280,1267.3800000000001,1282.83,1166.3599999999999,968.75,896.24,903.25,825.88,735.09,797.87,872.81,919.14,919.32,987.48,1020.62,1057.08,1036.19,1095.6300000000001,1115.0999999999999,1073.8699999999999,1104.49,1169.43,1186.69,1089.4100000000001,1030.71,1101.5999999999999,1049.33,1141.2,1183.26,1180.55,1257.6400000000001,1286.1199999999999,1327.22,1325.83,1363.61,1345.2,1320.64,1292.28,1218.8900000000001,1131.42,1253.3,1246.96,1257.5999999999999,1312.41,1365.68,1408.47,1397.91,1310.33,1362.16,1379.32,1406.58,1440.67,1412.16,1416.18,1426.19,1498.11,1514.68,1569.19,1597.57,1630.74,1606.28,1685.73,1632.97,1681.55,1756.54,1805.81,1848.36,1782.59,1859.45,1872.34,1883.95,1923.57,1960.23,1930.67,2003.37,1972.29,2018.05,2067.56,2058.9,1994.99,2104.5,2067.89,2085.5100000000002,2107.39,2063.11,2103.84,1972.18,1920.03,2079.36,2080.41,2043.94,1940.24,1932.23,2059.7399999999998,2065.3000000000002,2096.96,2098.86,2173.6,2170.9499999999998,2168.27,2126.15,2198.81,2238.83,2278.87,2363.64,2362.7199999999998,2384.1999999999998,2411.8000000000002,2423.41,2470.3000000000002"

$bParts = $bCsv -split ","
$rParts = $rCsv -split ","

$n = $bParts.Count

$bArr = New-Object 'object[,]' $n,1
$rArr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $bArr[$i,0] = [double]$bParts[$i]
    $rArr[$i,0] = [double]$rParts[$i]
}

$wsRaw.Range("B2:B206").Value = $bArr
$wsRaw.Range("R2:R206").Value = $rArr

# ---------------------------------------------------------------------------
# 2. Selections + active sheet/tab.
# ---------------------------------------------------------------------------
$wsData.Range("E10").Select()

$wsRaw.Activate()
$wsRaw.Range("D5").Select()

Write-Output "done"
